$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.685.72'

$ws.Range('D3').Value = '1.645.64'
$ws.Range('E3').Value = '  +1.17%  '

$ws.Range('E4').Value = '  +0.36%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '215.71'
$ws.Range('E5').Value = '  +1.18%  '

$ws.Range('E6').Value = '  +1.33%  '

$ws.Range('E7').Value = '  +0.22%  '

$ws.Range('E8').Value = '  +1.24%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0626'
$ws.Range('E9').Value = '  +0.49%  '

$ws.Range('E10').Value = '  +0.85%  '

$ws.Range('D12').Value = '1.873.09'
$ws.Range('E12').Value = '  +1.09%  '

$ws.Range('D13').Value = '1.629.32'
$ws.Range('E13').Value = '  -0.29%  '

$ws.Range('E14').Value = '  +0.95%  '

$ws.Range('E15').Value = '  +1.76%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.05'
$ws.Range('E16').Value = '  +0.60%  '

$ws.Range('D17').Value = '26.689.93'
$ws.Range('E17').Value = '  +0.91%  '

$ws.Range('D18').Value = '0.0₃0745'
$ws.Range('E18').Value = '  +0.62%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '217.00'
$ws.Range('E19').Value = '  +0.93%  '

$ws.Range('E20').Value = '  +0.22%  '

$ws.Range('E21').Value = '  +1.10%  '

$ws.Range('E22').Value = '  +0.31%  '

$ws.Range('E23').Value = '  +2.15%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.27'
$ws.Range('E24').Value = '  +13.67%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '145.63'
$ws.Range('E25').Value = '  -1.40%  '

$ws.Range('E26').Value = '  +0.35%  '

$ws.Range('E27').Value = '  +0.35%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.12'
$ws.Range('E28').Value = '  +4.18%  '

$ws.Range('E29').Value = '  +1.06%  '

$ws.Range('E30').Value = '  +0.90%  '

$ws.Range('E31').Value = '  +1.37%  '

$ws.Range('E32').Value = '  +1.41%  '

$ws.Range('E33').Value = '  +1.99%  '

$ws.Range('D34').Value = '1.276.41'
$ws.Range('E34').Value = '  +4.85%  '

$ws.Range('E35').Value = '  +3.45%  '

$ws.Range('E36').Value = '  +1.55%  '

$ws.Range('E37').Value = '  +3.05%  '

$ws.Range('E38').Value = '  +6.10%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.822'
$ws.Range('E39').Value = '  +3.51%  '

$ws.Range('E40').Value = '  +0.24%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.813'
$ws.Range('E41').Value = '  +2.32%  '

$ws.Range('E42').Value = '  +0.10%  '

$ws.Range('E43').Value = '  +1.31%  '

$ws.Range('D44').Value = '1.783.10'
$ws.Range('E44').Value = '  +1.15%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '91.90'
$ws.Range('E45').Value = '  -1.12%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '59.54'
$ws.Range('E46').Value = '  +8.64%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0516'
$ws.Range('E48').Value = '  +1.20%  '

$ws.Range('E49').Value = '  +3.14%  '

$ws.Range('E50').Value = '  +1.48%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.407'
$ws.Range('E51').Value = '  -0.02%  '
